$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (the previous "last row") loses its date-only format and becomes a
# regular datetime-formatted row, matching the format used by rows 2-15.
$ws.Range("A16").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add the new daily row (row 17) with the next day's data. It becomes the
# new "last row" and gets the date-only number format.
$ws.Range("A17").Value = 45966
$ws.Range("A17").NumberFormat = "YYYY-MM-DD"
$ws.Range("B17").Value = 35
$ws.Range("C17").Value = 42
$ws.Range("D17").Value = 43
